# Swap the two theme color schemes that ship with this deck.
#
# Before the edit:
#   ppt/theme/theme1.xml  (used by the notes master)   = generic "Office Theme" colors
#   ppt/theme/theme2.xml  (used by the slide master /
#                          presentation default theme)  = "Integral" theme colors
#
# After the edit the *content* of the two theme parts is swapped, so the
# slide master/presentation now carries the generic "Office Theme" palette
# and the notes master carries the "Integral" palette. The relationships
# (which part is named theme1.xml/theme2.xml) are untouched - only the
# colors inside each part change.
#
# PowerPoint's object model exposes the modern 12-slot DrawingML color
# scheme through Slide.ThemeColorScheme (Dark1, Light1, Dark2, Light2,
# Accent1-6, Hyperlink, FollowedHyperlink, in that order) which maps onto
# the <a:clrScheme> that the active slide's master/theme uses. We drive
# that to re-color the theme behind the slides (theme2.xml) with the
# original "Office Theme" palette.

function ConvertTo-VbaRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation

# Target palette = the generic "Office Theme" colors that used to live in
# theme1.xml, in ThemeColorScheme index order:
# 1 Dark1, 2 Light1, 3 Dark2, 4 Light2, 5 Accent1, 6 Accent2, 7 Accent3,
# 8 Accent4, 9 Accent5, 10 Accent6, 11 Hyperlink, 12 FollowedHyperlink
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-VbaRgb $officeThemeColors[$i - 1]
}
